$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a brand-new weekly record at row 24; this pushes the existing
# rows 24..99 down to 25..100 (dimension grows from A1:R99 to A1:R100).
$ws.Rows.Item(24).Insert()

$ws.Cells.Item(24, 1).Value = 7
$ws.Cells.Item(24, 2).Value = "Terminal Hortofrutícola Agro Chillán"
$ws.Cells.Item(24, 3).Value = "Ñuble"
$ws.Cells.Item(24, 4).Value = 44659
$ws.Cells.Item(24, 5).Value = 16
$ws.Cells.Item(24, 6).Value = 100112030
$ws.Cells.Item(24, 7).Value = "Poroto granado"
$ws.Cells.Item(24, 8).Value = "Sin especificar"
$ws.Cells.Item(24, 9).Value = "Primera"
$ws.Cells.Item(24, 10).Value = 100
$ws.Cells.Item(24, 11).Value = 19000
$ws.Cells.Item(24, 12).Value = 20000
$ws.Cells.Item(24, 13).Value = 19500
$ws.Cells.Item(24, 14).Value = "`$/saco 25 kilos"
$ws.Cells.Item(24, 15).Value = "Provincia de Diguillín"
$ws.Cells.Item(24, 16).Value = 780
$ws.Cells.Item(24, 17).Value = 25
$ws.Cells.Item(24, 18).Value = "Hortaliza"
